$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - new data row, mirrors the existing row 5 layout/styles.
$ws.Range("A6").Value = "26/06/2023 10:34"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 0.3
$ws.Range("B6").NumberFormat = "0%"

$ws.Range("C6").Formula = "=D6 * 600 / 100 *100"
$ws.Range("C6").NumberFormat = '0\ "mL"'

$ws.Range("D6").Value = 0.84
$ws.Range("D6").NumberFormat = "0%"

$ws.Range("E6").Formula = "=B6 * 600 / 100 *100"
$ws.Range("E6").NumberFormat = '0\ "mL"'

$ws.Range("F6").Value = 0.44
$ws.Range("F6").NumberFormat = "0%"

$ws.Range("G6").Formula = "=F6 * 600 / 100 *100"
$ws.Range("G6").NumberFormat = '0\ "mL"'

$ws.Range("H6").Value = 0.2
$ws.Range("H6").NumberFormat = "0%"

$ws.Range("I6").Formula = "=H6 * 600 / 100 *100"
$ws.Range("I6").NumberFormat = '0\ "mL"'

$ws.Range("J6").Value = 0.71
$ws.Range("J6").NumberFormat = "0%"

$ws.Range("K6").Formula = "=J6 * 600 / 100 *100"
$ws.Range("K6").NumberFormat = '0\ "mL"'

$ws.Range("L6").Value = 0.01
$ws.Range("L6").NumberFormat = "0%"

$ws.Range("M6").Formula = "=L6 * 600 / 100 *100"
$ws.Range("M6").NumberFormat = '0\ "mL"'

$ws.Range("N6").Formula = "=C6+E6+G6+I6+K6+M6"
$ws.Range("N6").NumberFormat = '0\ "mL"'
